$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.691.98'
$ws.Range("E2").Value = '  -2.24%  '
$ws.Range("D3").Value = '2.367.96'
$ws.Range("E3").Value = '  -2.64%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.06'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.50'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.90%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.528'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '2.363.98'
$ws.Range("E9").Value = '  -2.26%  '
$ws.Range("E10").Value = '  -2.99%  '
$ws.Range("E11").Value = '  -1.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.08'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.17%  '
$ws.Range("E13").Value = '  -0.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.54'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.00%  '
$ws.Range("E16").Value = '  -3.13%  '
$ws.Range("D17").Value = '59.662.09'
$ws.Range("E17").Value = '  -2.24%  '
$ws.Range("D18").Value = '2.361.82'
$ws.Range("E18").Value = '  -2.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.00'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +9.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.45'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.58'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.01'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.73%  '
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("E25").Value = '  -3.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.05'
$ws.Range("D26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '559.39'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.60%  '
$ws.Range("E28").Value = '  -6.47%  '
$ws.Range("E29").Value = '  -2.32%  '
$ws.Range("D30").Value = '0.0₃0920'
$ws.Range("E30").Value = '  +1.22%  '
$ws.Range("E31").Value = '  +1.69%  '
$ws.Range("E32").Value = '  -3.12%  '
$ws.Range("E33").Value = '  -3.34%  '
$ws.Range("E34").Value = '  -1.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.41%  '
$ws.Range("E36").Value = '  +2.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '153.08'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.87%  '
$ws.Range("E38").Value = '  -0.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.52'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.14'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.97'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.80%  '
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.42'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.71%  '
$ws.Range("E44").Value = '  -1.36%  '
$ws.Range("E45").Value = '  +3.57%  '
$ws.Range("D46").Value = '0.0₆0299'
$ws.Range("E46").Value = '  +5.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '138.41'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.69%  '
$ws.Range("E48").Value = '  +0.29%  '
$ws.Range("E49").Value = '  -1.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0500'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.01'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.58%  '
